# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column D ("municipio-nombre") moves from being a "medida" (measure, xsd:int)
# to a "dim" (dimension) that uses the shared sdmx-dimension:refArea concept
# scheme and a URI-Municipio type.
#
# Column F ("aragon") stays a "dim" but switches from its own
# iaest-dimension:aragon concept scheme / skos:Concept / mapping-aragon.xlsx
# triple to the shared sdmx-dimension:refArea concept scheme and a
# URI-Comunidad type, no longer needing a dedicated mapping file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: municipio-nombre
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column F: aragon
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("F5").Clear()
